# Auto-generated PowerShell script to reorder "Recorded By" (column G) values
# by swapping the last two comma-separated entries, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -> "System, system, backup@backdoor.com"
$rows = @(2, 28, 54)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = "System, system, backup@backdoor.com"
}

# -> "admin@admin.com, System"
$rows = @(7, 33, 59)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = "admin@admin.com, System"
}

# -> "System, dnasr281@gmail.com"
$rows = @(10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 26, 36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 52, 62, 63, 64, 65, 66, 67, 69, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86, 90, 92, 93, 94, 96, 99, 101, 109, 110, 111, 112, 116, 118, 119, 120, 122, 125, 127, 135, 136, 137, 138, 142, 144, 145, 146, 148, 151, 153)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = "System, dnasr281@gmail.com"
}

# -> "admin@admin.com, dnasr281@gmail.com"
$rows = @(87, 113, 139)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = "admin@admin.com, dnasr281@gmail.com"
}
